$wb = $excel.ActiveWorkbook

# --- Efficiency_SOEC: fix electrolyzer efficiency table (column B) ---
# Column B ("Efficiency [%]") was stored as a whole-number percentage (e.g. 3.452
# meaning 3.452%) while column A ("Power [%]") is a fraction (e.g. 0.089 meaning
# 8.9%). Convert column B to the same fractional convention and apply a
# percentage number format so the displayed numbers stay the same.
$wsSOEC = $wb.Worksheets.Item("Efficiency_SOEC")

$rngB = $wsSOEC.Range("B2:B32")
$rngB.Style = "Normal"
$rngB.NumberFormat = "0.00%"

for ($r = 2; $r -le 32; $r++) {
    $cell = $wsSOEC.Cells.Item($r, 2)
    $oldValue = $cell.Value()
    $cell.Value = $oldValue / 100
}

# --- Update the remembered cell selection on a couple of sheets ---
$wsPEM = $wb.Worksheets.Item("Efficiency_PEM")
$wsPEM.Activate()
[void]$wsPEM.Range("H5").Select()

$wsSOEC.Activate()
[void]$wsSOEC.Range("F3").Select()
